# Apply the StructureDefinition "reference-period" update:
#  - Version bumped 5.0.0 -> 6.0.0
#  - Date bumped to the new publish timestamp
#  - Publisher filled in ("Alvearie Team")
#  - Contact/"No display for ContactDetail" row replaced by a new
#    Jurisdiction/"United States of America" row
#  - the duplicated "Contact" row removed
#  - the Elements sheet's root Extension row gets a real Short/Definition
#    ("Reference Period" / "A time period in which the reference is valid")
#    instead of the generic placeholder text

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Version
$meta.Range("B3").Value = "6.0.0"

# Date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher (was blank)
$meta.Range("B9").Value = "Alvearie Team"

# Remove the second, duplicated "Contact" row (row 11); this shifts
# "Description" and everything below it up by one row.
$meta.Rows.Item(11).Delete()

# The remaining "Contact" row (now row 10) becomes the new
# Jurisdiction / United States of America row.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: give it a real short description / definition.
$elements.Range("K2").Value = "Reference Period"
$elements.Range("L2").Value = "A time period in which the reference is valid"
